$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A67").Value = "GRT-USD"
